$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Convert the phone-number cells from numeric values to text (shared strings)
$ws.Range("D1").Value = "714-943-9606"
$ws.Range("D2").Value = "714-555-5555"
$ws.Range("D3").Value = "555-555-5555"
$ws.Range("D4").Value = "714-444-4444"

# Update the active cell selection on the sheet
$ws.Range("E8").Select()
